# Apply updated TPM-derived values to Il23a-Il12rb1 LR-pair sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "G2" = 3.109891
    "H2" = 9.329673
    "I2" = 0.4376734453228661
    "J2" = 0.4376734453228661
    "K2" = 3
    "L2" = 1
    "M2" = 0.5699070000000001
    "N2" = 1.709721
    "O2" = 0.1475417602688563
    "P2" = 0.1475417602688563
    "Q2" = 1.772348650137
    "R2" = 15.951137851233
    "S2" = 0.06457511054587071
    "T2" = 0.0645751105458707
    "G3" = 3.109891
    "H3" = 9.329673
    "I3" = 0.4376734453228661
    "J3" = 0.4376734453228661
    "N3" = 5.084895
    "O3" = 0.4388051378454766
    "P3" = 0.4388051378454766
    "Q3" = 5.271156398815
    "R3" = 47.440407589335
    "S3" = 0.1920533565062049
    "T3" = 0.1920533565062049
    "G4" = 3.109891
    "H4" = 9.329673
    "I4" = 0.4376734453228661
    "J4" = 0.4376734453228661
    "M4" = 0.125128
    "N4" = 0.375384
    "O4" = 0.03239406671425592
    "P4" = 0.03239406671425593
    "Q4" = 0.3891344410479999
    "R4" = 3.502209969432
    "S4" = 0.01417802278684717
    "T4" = 0.01417802278684717
    "G5" = 3.109891
    "H5" = 9.329673
    "I5" = 0.4376734453228661
    "J5" = 0.4376734453228661
    "M5" = 1.472682666666667
    "N5" = 4.418048
    "O5" = 0.3812590351714111
    "P5" = 0.3812590351714111
    "Q5" = 4.579882570922666
    "R5" = 41.218943138304
    "S5" = 0.1668669554839433
    "T5" = 0.1668669554839433
    "I6" = 0.417132940486147
    "J6" = 0.417132940486147
    "K6" = 3
    "L6" = 1
    "M6" = 0.5699070000000001
    "N6" = 1.709721
    "O6" = 0.1475417602688563
    "P6" = 0.1475417602688563
    "Q6" = 1.689170343549
    "R6" = 15.202533091941
    "S6" = 0.06154452830545023
    "T6" = 0.06154452830545021
    "I7" = 0.417132940486147
    "J7" = 0.417132940486147
    "N7" = 5.084895
    "O7" = 0.4388051378454766
    "P7" = 0.4388051378454766
    "R7" = 45.213976143795
    "S7" = 0.1830400774499127
    "T7" = 0.1830400774499127
    "I8" = 0.417132940486147
    "J8" = 0.417132940486147
    "M8" = 0.125128
    "N8" = 0.375384
    "O8" = 0.03239406671425592
    "P8" = 0.03239406671425593
    "Q8" = 0.3708719260293333
    "R8" = 3.337847334264
    "S8" = 0.01351263230282199
    "T8" = 0.013512632302822
    "I9" = 0.417132940486147
    "J9" = 0.417132940486147
    "M9" = 1.472682666666667
    "N9" = 4.418048
    "O9" = 0.3812590351714111
    "P9" = 0.3812590351714111
    "Q9" = 4.364943553934222
    "R9" = 39.284491985408
    "S9" = 0.1590357024279621
    "T9" = 0.1590357024279621
    "G10" = 0.8168863333333333
    "H10" = 2.450659
    "I10" = 0.1149652691837634
    "J10" = 0.1149652691837634
    "K10" = 3
    "L10" = 1
    "M10" = 0.5699070000000001
    "N10" = 1.709721
    "O10" = 0.1475417602688563
    "P10" = 0.1475417602688563
    "Q10" = 0.465549239571
    "R10" = 4.189943156139
    "S10" = 0.01696217818515536
    "T10" = 0.01696217818515536
    "G11" = 0.8168863333333333
    "H11" = 2.450659
    "I11" = 0.1149652691837634
    "J11" = 0.1149652691837634
    "N11" = 5.084895
    "O11" = 0.4388051378454766
    "P11" = 0.4388051378454766
    "Q11" = 1.384593743978333
    "R11" = 12.461343695805
    "S11" = 0.05044735079162364
    "T11" = 0.05044735079162364
    "G12" = 0.8168863333333333
    "H12" = 2.450659
    "I12" = 0.1149652691837634
    "J12" = 0.1149652691837634
    "M12" = 0.125128
    "N12" = 0.375384
    "O12" = 0.03239406671425592
    "P12" = 0.03239406671425593
    "Q12" = 0.1022153531173333
    "R12" = 0.919938178056
    "S12" = 0.003724192599761223
    "T12" = 0.003724192599761224
    "G13" = 0.8168863333333333
    "H13" = 2.450659
    "I13" = 0.1149652691837634
    "J13" = 0.1149652691837634
    "M13" = 1.472682666666667
    "N13" = 4.418048
    "O13" = 0.3812590351714111
    "P13" = 0.3812590351714111
    "Q13" = 1.203014343736889
    "R13" = 10.827129093632
    "S13" = 0.0438315476072232
    "T13" = 0.0438315476072232
    "G14" = 0.2147876666666667
    "H14" = 0.644363
    "I14" = 0.03022834500722351
    "J14" = 0.03022834500722351
    "K14" = 3
    "L14" = 1
    "M14" = 0.5699070000000001
    "N14" = 1.709721
    "O14" = 0.1475417602688563
    "P14" = 0.1475417602688563
    "Q14" = 0.122408994747
    "R14" = 1.101680952723
    "S14" = 0.004459943232380052
    "T14" = 0.00445994323238005
    "G15" = 0.2147876666666667
    "H15" = 0.644363
    "I15" = 0.03022834500722351
    "J15" = 0.03022834500722351
    "N15" = 5.084895
    "O15" = 0.4388051378454766
    "P15" = 0.4388051378454766
    "Q15" = 0.3640575774316667
    "R15" = 3.276518196885001
    "S15" = 0.01326435309773534
    "T15" = 0.01326435309773534
    "G16" = 0.2147876666666667
    "H16" = 0.644363
    "I16" = 0.03022834500722351
    "J16" = 0.03022834500722351
    "M16" = 0.125128
    "N16" = 0.375384
    "O16" = 0.03239406671425592
    "P16" = 0.03239406671425593
    "Q16" = 0.02687595115466667
    "R16" = 0.241883560392
    "S16" = 0.0009792190248255435
    "T16" = 0.0009792190248255435
    "G17" = 0.2147876666666667
    "H17" = 0.644363
    "I17" = 0.03022834500722351
    "J17" = 0.03022834500722351
    "M17" = 1.472682666666667
    "N17" = 4.418048
    "O17" = 0.3812590351714111
    "P17" = 0.3812590351714111
    "Q17" = 0.3163140737137778
    "R17" = 2.846826663424
    "S17" = 0.01152482965228258
    "T17" = 0.01152482965228258
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
